$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) is treated as text so values like "1.001" are not
# auto-converted to numbers by Excel. Reset the style afterwards so no stray
# style index is left on the cells (keeps cells styleless, matching the source).
$priceRange = $ws.Range("D2:D51")
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '30.826.75'
$ws.Range("E2").Value = '  -0.93%  '
$ws.Range("D3").Value = '1.940.72'
$ws.Range("E3").Value = '  -0.64%  '
$ws.Range("D4").Value = '1.001'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '242.93'
$ws.Range("E5").Value = '  -1.27%  '
$ws.Range("E6").Value = '  +0.02%  '
$ws.Range("D7").Value = '0.4894'
$ws.Range("E7").Value = '  +0.39%  '
$ws.Range("D8").Value = '0.2944'
$ws.Range("E8").Value = '  -0.68%  '
$ws.Range("D9").Value = '0.06916'
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = '19.49'
$ws.Range("E10").Value = '  +2.38%  '
$ws.Range("D11").Value = '106.06'
$ws.Range("E11").Value = '  +0.17%  '
$ws.Range("D12").Value = '1.947.11'
$ws.Range("E12").Value = '  -0.21%  '
$ws.Range("D13").Value = '0.07725'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '5.364'
$ws.Range("E14").Value = '  -0.63%  '
$ws.Range("D15").Value = '0.6999'
$ws.Range("E15").Value = '  -1.64%  '
$ws.Range("E16").Value = '  -3.60%  '
$ws.Range("D17").Value = '30.837.16'
$ws.Range("E17").Value = '  -0.53%  '
$ws.Range("D18").Value = '0.000007734'
$ws.Range("E18").Value = '  -0.02%  '
$ws.Range("D19").Value = '13.12'
$ws.Range("E19").Value = '  -0.64%  '
$ws.Range("D20").Value = '2.215.37'
$ws.Range("E20").Value = '  +0.99%  '
$ws.Range("D21").Value = '1.000'
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").Value = '5.517'
$ws.Range("E22").Value = '  -0.06%  '
$ws.Range("D23").Value = '1.001'
$ws.Range("E23").Value = '  +0.05%  '
$ws.Range("D24").Value = '6.580'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").Value = '9.739'
$ws.Range("E25").Value = '  -1.61%  '
$ws.Range("D26").Value = '166.93'
$ws.Range("E26").Value = '  -0.85%  '
$ws.Range("D27").Value = '19.63'
$ws.Range("E27").Value = '  -1.23%  '
$ws.Range("D28").Value = '2.169'
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").Value = '0.1041'
$ws.Range("E29").Value = '  -0.91%  '
$ws.Range("E30").Value = '  -3.36%  '
$ws.Range("D31").Value = '4.580'
$ws.Range("E31").Value = '  -3.05%  '
$ws.Range("D32").Value = '1.555'
$ws.Range("E32").Value = '  -2.37%  '
$ws.Range("D33").Value = '4.374'
$ws.Range("E33").Value = '  -2.06%  '
$ws.Range("D34").Value = '0.04866'
$ws.Range("E34").Value = '  -2.55%  '
$ws.Range("D35").Value = '0.7567'
$ws.Range("E35").Value = '  -0.42%  '
$ws.Range("D36").Value = '1.158'
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("D37").Value = '0.9996'
$ws.Range("E37").Value = '  +0.01%  '
$ws.Range("E38").Value = '  +0.06%  '
$ws.Range("E39").Value = '  -1.95%  '
$ws.Range("B40").Value = 'MXToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D40").Value = '2.660'
$ws.Range("E40").Value = '  -2.01%  '
$ws.Range("B41").Value = 'FraxShare'
$ws.Range("C41").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D41").Value = '6.538'
$ws.Range("E41").Value = '  +1.83%  '
$ws.Range("D42").Value = '77.60'
$ws.Range("E42").Value = '  +6.77%  '
$ws.Range("D43").Value = '2.098'
$ws.Range("E43").Value = '  -2.12%  '
$ws.Range("D44").Value = '0.9065'
$ws.Range("E44").Value = '  +2.87%  '
$ws.Range("E45").Value = '  -1.24%  '
$ws.Range("D46").Value = '108.16'
$ws.Range("E46").Value = '  -1.29%  '
$ws.Range("D47").Value = '0.9988'
$ws.Range("E47").Value = '  -0.05%  '
$ws.Range("D48").Value = '7.780'
$ws.Range("E48").Value = '  +4.21%  '
$ws.Range("D49").Value = '995.59'
$ws.Range("E49").Value = '  +1.30%  '
$ws.Range("D50").Value = '0.1249'
$ws.Range("E50").Value = '  -2.10%  '
$ws.Range("B51").Value = 'EnergySwap'
$ws.Range("C51").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D51").Value = '9.313'
$ws.Range("E51").Value = '  -0.77%  '

# Restore default (unstyled) cell style on the Price column now that the text
# values are safely stored, so the cells keep no explicit style like the source.
$priceRange.Style = "Normal"
